$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new station rows (NAA stations, buffers) into the table.
# Row 35: RS107 / Sellingen
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "RS107"
$ws.Range("B35").Value = "Sellingen"

# Row 37: RS206 / Het Haantje
$ws.Rows.Item(37).Insert()
$ws.Range("A37").Value = "RS206"
$ws.Range("B37").Value = "Het Haantje"

# Row 38: RS207 / Nieuw-Dordrecht
$ws.Rows.Item(38).Insert()
$ws.Range("A38").Value = "RS207"
$ws.Range("B38").Value = "Nieuw-Dordrecht"

# Row 48: RS411 / Taarlo
$ws.Rows.Item(48).Insert()
$ws.Range("A48").Value = "RS411"
$ws.Range("B48").Value = "Taarlo"

# Restore the selection/view state recorded in the saved workbook.
$ws.Range("C44").Select() | Out-Null
